# Update "想去人数" (want-to-go count, column F) figures across the
# workbook's sheets to match the freshly scraped snapshot
# (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 1088
$ws.Range("F13").Value = 1408
$ws.Range("F15").Value = 1426
$ws.Range("F17").Value = 1174
$ws.Range("F18").Value = 284
$ws.Range("F20").Value = 439
$ws.Range("F21").Value = 770
$ws.Range("F23").Value = 47
$ws.Range("F25").Value = 1276
$ws.Range("F30").Value = 1029
$ws.Range("F36").Value = 1041
$ws.Range("F37").Value = 26
$ws.Range("F39").Value = 1564

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 122
$ws.Range("F34").Value = 56
$ws.Range("F46").Value = 117
$ws.Range("F47").Value = 53

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F12").Value = 177
$ws.Range("F13").Value = 667
$ws.Range("F15").Value = 352

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value  = 177
$ws.Range("F9").Value  = 177
$ws.Range("F10").Value = 667
$ws.Range("F11").Value = 667
$ws.Range("F14").Value = 1088
$ws.Range("F24").Value = 1408
$ws.Range("F26").Value = 1426
$ws.Range("F27").Value = 1174
$ws.Range("F31").Value = 770
$ws.Range("F33").Value = 352
$ws.Range("F35").Value = 1276
$ws.Range("F38").Value = 1029
$ws.Range("F42").Value = 26
$ws.Range("F45").Value = 1564
